# Revert "Merge branch 'wrong-xlsform-col'"
#
# The "survey" sheet's header row (row 1) had its C column label cell
# accidentally changed to "label" by the bad merge; reverting restores it
# to "message" (the XLSForm "survey" sheet's third header column is the
# validation/constraint "message" column, not the "label" column used on
# the "choices" sheet).
#
# choices!C1 legitimately stays "label" (XLSForm choices sheets use
# type/name/label) - only the survey sheet's header is wrong and needs
# reverting.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item(1)
$survey.Range("C1").Value() = "message"

# Cosmetic: restore the previously-selected cell on the survey sheet.
$null = $survey.Range("C3").Select()
